# The post at row 383 ("恥ずかしがり屋のマヌルネコと暖かい朝の太陽光線") was removed.
# Deleting the entire row shifts every subsequent row up by one and
# automatically shrinks the sheet's used range from A1:C415 to A1:C414,
# matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(383).Delete()
